$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-27 Saturday", "2025-09-28 Sunday"),
    @("31÷9=", "81÷3="),
    @("22÷2=", "19÷5="),
    @("18÷3=", "65÷4="),
    @("27÷2=", "93÷8="),
    @("88÷3=", "79÷8="),
    @("76÷8=", "76÷7="),
    @("41÷9=", "98÷4="),
    @("68÷7=", "12÷6="),
    @("38÷9=", "93÷3="),
    @("16÷3=", "47÷4="),
    @("69÷8=", "49÷7="),
    @("58÷9=", "29÷2="),
    @("19÷4=", "75÷8="),
    @("38÷4=", "39÷8="),
    @("93÷5=", "28÷7="),
    @("30÷5=", "41÷7="),
    @("16÷6=", "29÷9="),
    @("47÷5=", "88÷2="),
    @("81÷8=", "67÷4="),
    @("42÷4=", "44÷3="),
    @("84÷8=", "69÷2="),
    @("34÷8=", "55÷5="),
    @("17÷8=", "16÷4="),
    @("55÷7=", "87÷9="),
    @("29÷7=", "86÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
